# Update the Mapping sheet bounding-box coordinates after converting all
# shapefiles to WGS 84 (EPSG 4269).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -82.5275
$ws.Range("B2").Value = -82.5262

$ws.Range("A3").Value = 41.3631
$ws.Range("B3").Value = 41.364

$ws.Range("A4").Value = -82.4981
$ws.Range("B4").Value = -82.4994

$ws.Range("A5").Value = 41.3871
$ws.Range("B5").Value = 41.3862

# Touch B6 so the cell becomes present (matches the source edit, which
# materializes an empty B6 cell in the sheet's XML).
$ws.Cells.Item(6, 2).Font.Bold = $false
$ws.Cells.Item(6, 2).ClearFormats()
